# AfDD_2022_Annex_Table_Tab17.xlsx - apply the tracked edit
#
# Summary of the change (per the OOXML diff):
#  1. Swap the "Exports of goods and services" and "Foreign balance"
#     headers/columns (L <-> N), for both the header row (row 2) and every
#     data row (rows 3-99). Column M ("Imports...") is untouched.
#  2. Introduce a typo in the footnote text: "explore" -> "Implore" (cell
#     B108), while the hyperlink attached to that cell keeps displaying the
#     original ("explore") wording as its stored display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header swap (row 2): L2 <-> N2 text ("Exports..." <-> "Foreign...")
# ---------------------------------------------------------------------
$headerL = $ws.Range("L2").Value2
$headerN = $ws.Range("N2").Value2
$ws.Range("L2").Value2 = $headerN
$ws.Range("N2").Value2 = $headerL

# ---------------------------------------------------------------------
# 2. Data swap (rows 3-99): L <-> N values, row by row
# ---------------------------------------------------------------------
for ($r = 3; $r -le 99; $r++) {
    $lCell = $ws.Range("L$r")
    $nCell = $ws.Range("N$r")
    $lVal = $lCell.Value2
    $nVal = $nCell.Value2
    if ($lVal -ne $nVal) {
        $lCell.Value2 = $nVal
        $nCell.Value2 = $lVal
    }
}

# ---------------------------------------------------------------------
# 3. Footnote typo: "explore" -> "Implore" in B108, but keep the
#    hyperlink's stored display text as the original wording.
# ---------------------------------------------------------------------
$originalText = "If you would like to explore these data further, look up historic values for these indicators, or produce interactive visualisations of these data, please visit the website https://oe.cd/AFDD-2022"
$typoText = "If you would like to Implore these data further, look up historic values for these indicators, or produce interactive visualisations of these data, please visit the website https://oe.cd/AFDD-2022"

# Re-assert the hyperlink on B108 with its original target, storing the
# original wording as the hyperlink's display text.
$b108 = $ws.Range("B108")
$ws.Hyperlinks.Add($b108, "https://oe.cd/AFDD-2022", "", "", $originalText) | Out-Null

# Now set the actual cell content to the typo'd wording without disturbing
# the hyperlink relationship that's already attached to B108.
$ws.Range("B108").Value = $typoText
